# Apply the updated cryptos list (Price / Volume(1h) columns)
# Values are written as text (matching the original inline-string cells),
# forcing a Text number format while the value is assigned so that
# numeric-looking strings (e.g. "0.548", "4.00") are not silently
# converted into floating point numbers, then resetting the style back
# to "Normal" so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# row => (new Price value or $null when unchanged, new Volume(1h) value)
$updates = @{
    2 = @("41.594.83", "  +0.27%  ")
    3 = @("2.465.15", "  +0.31%  ")
    4 = @($null, "  -0.63%  ")
    5 = @("314.56", "  +1.11%  ")
    6 = @("91.25", "  +1.12%  ")
    7 = @("0.548", "  +2.37%  ")
    8 = @($null, "  -0.82%  ")
    9 = @("0.512", "  +5.51%  ")
    10 = @("32.56", "  +0.35%  ")
    11 = @($null, "  +3.11%  ")
    12 = @($null, "  +0.72%  ")
    13 = @("2.845.03", "  +0.10%  ")
    14 = @("6.88", "  +1.80%  ")
    15 = @("15.77", "  +4.80%  ")
    16 = @("2.459.41", "  -2.77%  ")
    17 = @("0.775", "  +1.03%  ")
    18 = @("41.540.68", "  +0.63%  ")
    19 = @("6.48", "  +4.49%  ")
    20 = @($null, "  +3.67%  ")
    21 = @("70.94", "  +1.49%  ")
    22 = @("11.34", "  +4.44%  ")
    23 = @("237.89", "  +1.70%  ")
    24 = @("2.71", "  +1.11%  ")
    25 = @("1.91", "  +3.18%  ")
    26 = @($null, "  -0.05%  ")
    27 = @("24.27", "  +2.17%  ")
    28 = @($null, "  +0.82%  ")
    29 = @("9.66", "  +0.93%  ")
    30 = @("35.28", "  -0.17%  ")
    31 = @("156.30", "  +3.00%  ")
    32 = @("5.44", "  +1.84%  ")
    33 = @($null, "  +0.84%  ")
    34 = @("0.0757", "  +1.75%  ")
    35 = @("17.22", "  +0.57%  ")
    36 = @($null, "  -8.50%  ")
    37 = @($null, "  -2.61%  ")
    38 = @($null, "  +2.51%  ")
    39 = @("0.103", "  +3.95%  ")
    40 = @($null, "  -2.06%  ")
    41 = @("4.00", "  +1.05%  ")
    42 = @($null, "  -1.29%  ")
    43 = @("1.958.66", "  +0.64%  ")
    44 = @($null, "  +1.91%  ")
    45 = @("18.73", "  -1.78%  ")
    46 = @($null, "  +0.73%  ")
    47 = @("9.00", "  +5.23%  ")
    48 = @("2.703.02", "  -0.48%  ")
    49 = @("96.78", "  +2.04%  ")
    50 = @("67.00", "  +0.66%  ")
    51 = @("0.172", "  -0.54%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $dVal = $pair[0]
    $eVal = $pair[1]
    if ($dVal -ne $null) {
        Set-TextValue $ws.Cells.Item($row, 4) $dVal
    }
    Set-TextValue $ws.Cells.Item($row, 5) $eVal
}
